$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (Anmeldung / Email input) : add a new reviewer comment in column E ---
$ws.Range("E11").Value = "Ist hier nur von der Eingabe die Rede oder auch von der Funktionalität? Dann weniger Punkte…"
$ws.Range("E11").Font.Italic = $true
$ws.Range("E11").VerticalAlignment = -4108

# --- Row 13 : update existing reviewer comment text ---
$ws.Range("E13").Value = "Als Banner?`nWas heißt ""dauerhaft sichtbar""?"

# --- Row 35 (Überfällige Todos / overdue highlighting) : mark points awarded + comment ---
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = "Uhrzeit wird immer als AM gespiechert."
$ws.Range("E35").Font.Italic = $true
$ws.Range("E35").VerticalAlignment = -4108

# --- Row 36 : remove the old reviewer comment ("Fehlt da was?") ---
$ws.Range("E36").Clear()

# --- Conditional formatting on column D: keep the overdue highlight pointed at the
#     correct (orange) fill used for mismatched grading cells ---
$fc = $ws.Range("D1:D1048576").FormatConditions
$cond = $fc.Item(1)
$cond.Interior.Color = 49407

# --- Update the view/selection to match where the review was focused ---
$ws.Range("E36").Select()
